$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked" and "is_enabled" dict columns (D1 / E1 previously held
# their header templates). Shift the remaining "order_by" / "rem" header
# templates into D1 / E1, and drop the now-trailing F1 / G1 cells entirely.
$ws.Range("D1").Value = "<%=comment.order_by%>"
$ws.Range("E1").Value = "<%=comment.rem%>"
$ws.Range("F1:G1").ClearContents()
